$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows to append (dates 2021-04-22 .. 2021-04-26 as Excel serials)
$newRows = @(
    @{ Row = 234; A = 44308; B = 1; C = 1; D = 53.73455131649651 },
    @{ Row = 235; A = 44309; B = 0; C = 1; D = 53.73455131649651 },
    @{ Row = 236; A = 44310; B = 0; C = 1; D = 53.73455131649651 },
    @{ Row = 237; A = 44311; B = 0; C = 1; D = 53.73455131649651 },
    @{ Row = 238; A = 44312; B = 0; C = 1; D = 53.73455131649651 }
)

foreach ($r in $newRows) {
    $rowNum = $r.Row

    # Column A: date value, styled like the existing date column (copy format from A233)
    $ws.Range("A233").Copy() | Out-Null
    $cellA = $ws.Range("A$rowNum")
    $cellA.PasteSpecial(-4122) | Out-Null   # xlPasteFormats
    $cellA.Value = $r.A

    $ws.Range("B$rowNum").Value = $r.B
    $ws.Range("C$rowNum").Value = $r.C
    $ws.Range("D$rowNum").Value = $r.D
}

$excel.CutCopyMode = 0
